# The deck's slide master (and therefore every slide) is wired to the
# "Integral" theme palette; the commit swaps the active theme's colour
# scheme over to the stock "Office Theme" palette (accent colours,
# dark/light pairs and hyperlink colours) that previously only lived,
# unused, in the Notes Master's theme part.
#
# PowerPoint's ThemeColorScheme.Colors(i).RGB indices map 1:1 onto the
# <a:clrScheme> slot order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# .RGB is an OLE COLORREF (0x00BBGGRR), so convert from the usual RRGGBB
# hex notation before assigning.

function ConvertTo-OleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return [System.Drawing.ColorTranslator]::ToOle([System.Drawing.Color]::FromArgb($r, $g, $b))
}

# Target palette ("Office Theme"), in clrScheme slot order.
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$p = $ppt.ActivePresentation
$themeColors = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-OleColor $officeThemeColors[$i - 1]
}
